$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT (matching the
# source inline-string cells) instead of letting Excel auto-convert it to a
# number. We do this by writing a text-formula ("=""123.45""") and then doing
# a Copy / Paste-Special-Values on that same cell, which bakes the formula
# down to a literal string value without touching the cell's style/numFmt.
function Set-TextNumber {
    param($cellRef, $val)
    $escaped = $val -replace '"', '""'
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

# --- Plain text / non-ambiguous values: assign directly ---
$ws.Range("D2").Value = '51.753.36'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.803.23'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +7.64%  '
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("E14").Value = '  +2.21%  '
$ws.Range("D15").Value = '3.244.77'
$ws.Range("E15").Value = '  +0.76%  '
$ws.Range("D16").Value = '2.817.86'
$ws.Range("E16").Value = '  +1.63%  '
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").Value = '51.760.05'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("E20").Value = '  +3.69%  '
$ws.Range("E21").Value = '  +3.24%  '
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("E28").Value = '  -4.46%  '
$ws.Range("E29").Value = '  +11.37%  '
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("E34").Value = '  +8.21%  '
$ws.Range("E35").Value = '  +6.01%  '
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  +2.36%  '
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("E42").Value = '  -1.31%  '
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("E46").Value = '  +4.53%  '
$ws.Range("D47").Value = '2.121.55'
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("E48").Value = '  +6.51%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("E51").Value = '  +7.08%  '

# --- Values that look like plain numbers: force text storage ---
Set-TextNumber "D5" "354.95"
Set-TextNumber "D6" "111.92"
Set-TextNumber "D9" "0.634"
Set-TextNumber "D10" "40.38"
Set-TextNumber "D13" "20.02"
Set-TextNumber "D17" "0.946"
Set-TextNumber "D19" "7.66"
Set-TextNumber "D21" "13.66"
Set-TextNumber "D23" "70.53"
Set-TextNumber "D24" "268.67"
Set-TextNumber "D27" "26.18"
Set-TextNumber "D28" "0.161"
Set-TextNumber "D29" "39.12"
Set-TextNumber "D30" "10.38"
Set-TextNumber "D32" "6.14"
Set-TextNumber "D33" "52.22"
Set-TextNumber "D35" "0.0885"
Set-TextNumber "D36" "0.0444"
Set-TextNumber "D38" "18.91"
Set-TextNumber "D42" "2.51"
Set-TextNumber "D45" "22.05"
Set-TextNumber "D46" "3.42"
Set-TextNumber "D48" "2.42"
Set-TextNumber "D49" "0.949"

$excel.CutCopyMode = 0

